$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.684.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.67%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.279.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.68%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.30%  '

$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.600'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0901'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.87%  '

$ws.Range("E12").Value = '  -4.94%  '

$ws.Range("E13").Value = '  -0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.956'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.624.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.280.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.813.96'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.30%  '

$ws.Range("E20").Value = '  -1.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '278.90'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.71%  '

$ws.Range("E25").Value = '  -3.72%  '

$ws.Range("E26").Value = '  +0.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.01%  '

$ws.Range("E28").Value = '  +3.63%  '

$ws.Range("E29").Value = '  -0.40%  '

$ws.Range("E30").Value = '  -5.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0866'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.25%  '

$ws.Range("E33").Value = '  -1.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.132'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("E36").Value = '  -6.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0344'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.71%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.62'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +15.83%  '

$ws.Range("E42").Value = '  -4.52%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.45%  '

$ws.Range("E44").Value = '  +0.21%  '

$ws.Range("E45").Value = '  -8.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '115.48'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.85%  '

$ws.Range("E51").Value = '  -4.34%  '
